$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 13891389
$ws.Range("J43").Value = 18521186
$ws.Range("L43").Value = 18521186
$ws.Range("N43").Value = -18521324
$ws.Range("H53").Value = 3113
$ws.Range("I53").Value = 2875.25
$ws.Range("J53").Value = 3430
$ws.Range("K53").Value = 2875.25
$ws.Range("L53").Value = 3430
$ws.Range("M53").Value = -2238.25
$ws.Range("N53").Value = -4704
$ws.Range("H74").Value = 9654.467000000001
$ws.Range("I74").Value = 10810.272
$ws.Range("J74").Value = 6476
$ws.Range("K74").Value = 10810.272
$ws.Range("L74").Value = 6476
$ws.Range("M74").Value = -9874.272000000001
$ws.Range("N74").Value = -8348
$ws.Range("H77").Value = 9654.467000000001
$ws.Range("I77").Value = 10810.272
$ws.Range("J77").Value = 6476
$ws.Range("K77").Value = 54051.36
$ws.Range("L77").Value = 32380
$ws.Range("M77").Value = -49371.36
$ws.Range("N77").Value = -41740
$ws.Range("H106").Value = 2801.2856
$ws.Range("I106").Value = 2921.2
$ws.Range("J106").Value = 2501.5
$ws.Range("K106").Value = 2921.2
$ws.Range("L106").Value = 2501.5
$ws.Range("M106").Value = -2290.2
$ws.Range("N106").Value = -3763.5
$ws.Range("H113").Value = 20003410
$ws.Range("I113").Value = 66668932
$ws.Range("J113").Value = 3899.7144
$ws.Range("K113").Value = 66668932
$ws.Range("L113").Value = 3899.7144
$ws.Range("M113").Value = -66665678
$ws.Range("N113").Value = -10407.7144
$ws.Range("H129").Value = 836.2879
$ws.Range("I129").Value = 404.8
$ws.Range("K129").Value = 1214.4
$ws.Range("M129").Value = 3785.6
$ws.Range("H132").Value = 7942456.5
$ws.Range("I132").Value = 10103508
$ws.Range("K132").Value = 30310524
$ws.Range("M132").Value = -30307994
$ws.Range("H137").Value = 1122.4584
$ws.Range("I137").Value = 1135.3334
$ws.Range("K137").Value = 3406.0002
$ws.Range("M137").Value = -856.0001999999999
$ws.Range("H138").Value = 506842.06
$ws.Range("I138").Value = 692.2778
$ws.Range("J138").Value = 871269.9399999999
$ws.Range("K138").Value = 2076.8334
$ws.Range("L138").Value = 2613809.82
$ws.Range("M138").Value = 3063.1666
$ws.Range("N138").Value = -2624089.82

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1701.3334
$ws.Range("I61").Value = 1187.4286
$ws.Range("K61").Value = 1187.4286
$ws.Range("M61").Value = -975.4286
$ws.Range("H63").Value = 2414.2856
$ws.Range("I63").Value = 2416.6667
$ws.Range("K63").Value = 2416.6667
$ws.Range("M63").Value = -1730.6667
$ws.Range("H66").Value = 2414.2856
$ws.Range("I66").Value = 2416.6667
$ws.Range("K66").Value = 12083.3335
$ws.Range("M66").Value = -8651.333500000001
$ws.Range("H74").Value = 1295.9166
$ws.Range("I74").Value = 873.2857
$ws.Range("K74").Value = 873.2857
$ws.Range("M74").Value = 0.7142999999999802
$ws.Range("H77").Value = 1295.9166
$ws.Range("I77").Value = 873.2857
$ws.Range("K77").Value = 4366.4285
$ws.Range("M77").Value = 1.571500000000015
$ws.Range("H132").Value = 1881
$ws.Range("I132").Value = 1507.4242
$ws.Range("J132").Value = 3113.8
$ws.Range("K132").Value = 4522.2726
$ws.Range("L132").Value = 9341.400000000001
$ws.Range("M132").Value = -1992.2726
$ws.Range("N132").Value = -14401.4
$ws.Range("H136").Value = 1701.3334
$ws.Range("I136").Value = 1187.4286
$ws.Range("K136").Value = 3562.2858
$ws.Range("M136").Value = -1012.2858

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5736.4287
$ws.Range("I86").Value = 5734.1665
$ws.Range("J86").Value = 5750
$ws.Range("K86").Value = 5734.1665
$ws.Range("L86").Value = 5750
$ws.Range("M86").Value = -4611.1665
$ws.Range("N86").Value = -7996
$ws.Range("H89").Value = 5736.4287
$ws.Range("I89").Value = 5734.1665
$ws.Range("J89").Value = 5750
$ws.Range("K89").Value = 28670.8325
$ws.Range("L89").Value = 28750
$ws.Range("M89").Value = -23054.8325
$ws.Range("N89").Value = -39982

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 15011.857
$ws.Range("I50").Value = 3041.5
$ws.Range("J50").Value = 19800
$ws.Range("K50").Value = 3041.5
$ws.Range("L50").Value = 19800
$ws.Range("M50").Value = -2416.5
$ws.Range("N50").Value = -21050
$ws.Range("H99").Value = 1463299.1
$ws.Range("I99").Value = 2632992.2
$ws.Range("K99").Value = 2632992.2
$ws.Range("M99").Value = -2631494.2
$ws.Range("H107").Value = 520.05554
$ws.Range("I107").Value = 373.81818
$ws.Range("J107").Value = 749.8570999999999
$ws.Range("K107").Value = 373.81818
$ws.Range("L107").Value = 749.8570999999999
$ws.Range("M107").Value = 1546.18182
$ws.Range("N107").Value = -4589.8571
$ws.Range("H126").Value = 1463299.1
$ws.Range("I126").Value = 2632992.2
$ws.Range("K126").Value = 7898976.600000001
$ws.Range("M126").Value = -7896506.600000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 595.1667
$ws.Range("I14").Value = 595.1667
$ws.Range("K14").Value = 1785.5001
$ws.Range("M14").Value = -1612.5001
$ws.Range("H47").Value = 203.16667
$ws.Range("I47").Value = 203.16667
$ws.Range("K47").Value = 609.50001
$ws.Range("M47").Value = -178.50001
$ws.Range("H131").Value = 37038760
$ws.Range("I131").Value = 166667000
$ws.Range("J131").Value = 2115.0952
$ws.Range("K131").Value = 500001000
$ws.Range("L131").Value = 6345.285600000001
$ws.Range("M131").Value = -499995960
$ws.Range("N131").Value = -16425.2856
$ws.Range("H139").Value = 3033
$ws.Range("I139").Value = 2812.25
$ws.Range("J139").Value = 3474.5
$ws.Range("K139").Value = 8436.75
$ws.Range("L139").Value = 10423.5
$ws.Range("M139").Value = -3296.75
$ws.Range("N139").Value = -20703.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 142.33333
$ws.Range("I2").Value = 75.75
$ws.Range("K2").Value = 75.75
$ws.Range("M2").Value = 37.25
$ws.Range("H70").Value = 28131800
$ws.Range("I70").Value = 27783602
$ws.Range("J70").Value = 28579486
$ws.Range("K70").Value = 27783602
$ws.Range("L70").Value = 28579486
$ws.Range("M70").Value = -27783332
$ws.Range("N70").Value = -28580026
$ws.Range("H73").Value = 28131800
$ws.Range("I73").Value = 27783602
$ws.Range("J73").Value = 28579486
$ws.Range("K73").Value = 27783602
$ws.Range("L73").Value = 28579486
$ws.Range("M73").Value = -27782666
$ws.Range("N73").Value = -28581358
$ws.Range("H122").Value = 1053.5
$ws.Range("I122").Value = 1053.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3160.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -710.5
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2372.2273
$ws.Range("I132").Value = 1993.1666
$ws.Range("J132").Value = 4078
$ws.Range("K132").Value = 5979.4998
$ws.Range("L132").Value = 12234
$ws.Range("M132").Value = -3449.4998
$ws.Range("N132").Value = -17294

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2662.9412
$ws.Range("I40").Value = 2413.077
$ws.Range("J40").Value = 3475
$ws.Range("K40").Value = 2413.077
$ws.Range("L40").Value = 3475
$ws.Range("M40").Value = -2277.077
$ws.Range("N40").Value = -3747
$ws.Range("H46").Value = 1637.2858
$ws.Range("I46").Value = 1451.6
$ws.Range("J46").Value = 2101.5
$ws.Range("K46").Value = 1451.6
$ws.Range("L46").Value = 2101.5
$ws.Range("M46").Value = -1263.6
$ws.Range("N46").Value = -2477.5
$ws.Range("H68").Value = 1662.1666
$ws.Range("I68").Value = 1593.8572
$ws.Range("J68").Value = 1901.25
$ws.Range("K68").Value = 1593.8572
$ws.Range("L68").Value = 1901.25
$ws.Range("M68").Value = -844.8571999999999
$ws.Range("N68").Value = -3399.25
$ws.Range("H71").Value = 1662.1666
$ws.Range("I71").Value = 1593.8572
$ws.Range("J71").Value = 1901.25
$ws.Range("K71").Value = 7969.286
$ws.Range("L71").Value = 9506.25
$ws.Range("M71").Value = -4225.286
$ws.Range("N71").Value = -16994.25
$ws.Range("H110").Value = 25933
$ws.Range("J110").Value = 23899.5
$ws.Range("L110").Value = 23899.5
$ws.Range("N110").Value = -32079.5
$ws.Range("H122").Value = 14168726
$ws.Range("I122").Value = 20239852
$ws.Range("K122").Value = 60719556
$ws.Range("M122").Value = -60717106

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 45462180
$ws.Range("I62").Value = 55562000
$ws.Range("J62").Value = 13000
$ws.Range("K62").Value = 55562000
$ws.Range("L62").Value = 13000
$ws.Range("M62").Value = -55561376
$ws.Range("N62").Value = -14248
$ws.Range("H65").Value = 45462180
$ws.Range("I65").Value = 55562000
$ws.Range("J65").Value = 13000
$ws.Range("K65").Value = 277810000
$ws.Range("L65").Value = 65000
$ws.Range("M65").Value = -277806880
$ws.Range("N65").Value = -71240
